$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded, which pushes the existing
# rows 46-91 down one row (to 47-92) and inserts a fresh row 46 with the
# new data (matches the Excel UI behaviour of inserting a row with
# xlShiftDown and then filling it in).
$ws.Rows(46).Insert()

$ws.Cells.Item(46, 1).Value = 7
$ws.Cells.Item(46, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(46, 3).Value = "Ñuble"
$ws.Cells.Item(46, 4).Value = 44904
$ws.Cells.Item(46, 5).Value = 16
$ws.Cells.Item(46, 6).Value = 100112022
$ws.Cells.Item(46, 7).Value = "Arveja Verde"
$ws.Cells.Item(46, 8).Value = "Sin especificar"
$ws.Cells.Item(46, 9).Value = "Primera"
$ws.Cells.Item(46, 10).Value = 80
$ws.Cells.Item(46, 11).Value = 22000
$ws.Cells.Item(46, 12).Value = 23000
$ws.Cells.Item(46, 13).Value = 22500
$ws.Cells.Item(46, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(46, 15).Value = "Región del Maule"
$ws.Cells.Item(46, 16).Value = 900
$ws.Cells.Item(46, 17).Value = 25
$ws.Cells.Item(46, 18).Value = "Hortaliza"
